# "hoan thanh add san pham" - populate the first product-detail row with
# its name, material, and weight.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("chi tiet sp")

# Pull the wrap-text / vertical-centered format that's already used on this
# sheet (L2/M2 style differs - it's the hyperlink style; the plain
# wrap+center style lives on "loai"!A4) and stamp it onto the three cells
# that are about to receive the new product data, without touching B2.
$xlPasteFormats = -4122
$fmtSrc = $wb.Worksheets.Item("loai").Range("A4")
$fmtSrc.Copy()
$ws.Range("A2").PasteSpecial($xlPasteFormats)
$ws.Range("C2").PasteSpecial($xlPasteFormats)
$ws.Range("D2").PasteSpecial($xlPasteFormats)

$ws.Range("A2").Value = "ROYAL M20D"
$ws.Range("C2").Value = "xốp EPS"
$ws.Range("D2").Value = 850

# Leave the selection where the author ended up after typing the new row.
$ws.Range("K8").Select() | Out-Null
